$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: labels - improve casing / punctuation
$ws.Range("A1").Value = "Horas trabajadas"
$ws.Range("B1").Value = "Personas residentes viviendas familiares"
$ws.Range("C1").Value = "Sector actividad, descripción"
$ws.Range("D1").Value = "Aragón"
$ws.Range("E1").Value = "Sector actividad, código"

# Row 4: datatypes - fix erroneous measure datatypes for horas-trabajadas and sector-actividad-descripcion
$ws.Range("A4").Value = "xsd:string"
$ws.Range("C4").Value = "xsd:string"
